$d = $word.ActiveDocument

# Locate the start paragraph ("Participation: We will:") and end paragraph
# (the very last paragraph in the document body) that bound the region we
# are replacing, then swap in the freshly authored OOXML for that whole span
# in one shot.
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    if ($startPara -eq $null -and $p.Range.Text.StartsWith("Participation: ")) {
        $startPara = $p
    }
}
$endPara = $d.Paragraphs.Last

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = '<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Participation</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> and Work allocation</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:t>We will:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>We will use the SCUM methodology to undertake out project execution</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Work will be allocated evenly during a group discussion at each sprint cycle</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Each have a fair amount of contribution towards each stage of the assignment</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Make sure all team members are not just participating but understanding the requirements.</w:t></w:r><w:r><w:tab/></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>A</w:t></w:r><w:r><w:t>ll</w:t></w:r><w:r><w:t xml:space="preserve"> members must</w:t></w:r><w:r><w:t xml:space="preserve"> attend every </w:t></w:r><w:r><w:t xml:space="preserve">meeting </w:t></w:r><w:r><w:t>and contribute to its conversations</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="17"/></w:numPr></w:pPr><w:r><w:t>Contribution must be equal in regards</w:t></w:r><w:r><w:t xml:space="preserve"> to the execution</w:t></w:r><w:r><w:t xml:space="preserve"> of the proj</w:t></w:r><w:r><w:t>ect</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Communication: </w:t></w:r><w:r><w:t>We will:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:t>Communicate via Email and G</w:t></w:r><w:r><w:t>it</w:t></w:r><w:r><w:t>Hub to keep the team informed on individual team member progress</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Ask for help </w:t></w:r><w:r><w:t>from other members if encountering problems</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="20"/></w:numPr></w:pPr><w:r><w:t>Document all progress with</w:t></w:r><w:r><w:t>in</w:t></w:r><w:r><w:t xml:space="preserve"> the group</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="720"/></w:pPr></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Meeting Guidelines: </w:t></w:r><w:r><w:t>We will:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr><w:r><w:t>Hold 1 S</w:t></w:r><w:r><w:t>kype meeting a week</w:t></w:r><w:r><w:t xml:space="preserve"> with timing agreed via email</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr><w:r><w:t>Include all members in meeting unless otherwise arranged</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr><w:r><w:t>Hold 1 physical stand-up after Wednesday lab</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="18"/></w:numPr></w:pPr><w:r><w:t>Document meeting in the wiki under blog post</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Problem Solving: </w:t></w:r><w:r><w:t>We will:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:t>Avoid placing blame when things go wrong</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Work together when an assigned piece of work is difficult to understand or </w:t></w:r><w:r><w:t xml:space="preserve">to </w:t></w:r><w:r><w:t>complete</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:t>Discuss the process regarding a problem and explore how it can be improved</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="19"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">An anonymous team vote will decide the outcome of </w:t></w:r><w:r><w:t>a</w:t></w:r><w:r><w:t xml:space="preserve"> dispute or indecision</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p/><w:p/><w:p/>'

$range.InsertXML($xml)

Write-Output "done"
